$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the shared string values (accented "Alliés" -> unaccented "Allies")
# while keeping the relationships in the sheet the same.
# C8 and C9 currently hold "Alliés"; set this first so "Allies" is appended
# to the shared string table before "Allies commerciaux".
$ws.Range("C8").Value = "Allies"
$ws.Range("C9").Value = "Allies"

# C4 and C5 currently hold "Alliés commerciaux"
$ws.Range("C4").Value = "Allies commerciaux"
$ws.Range("C5").Value = "Allies commerciaux"

# Update the active cell selection to C5
$ws.Range("C5").Select()
